# Weekly fruit/vegetable price update: insert 3 new rows of data
# (new "Sandia" / O'Higgins region prices dated 2022-02-18, serial 44610)
# above the existing row 93, pushing the old rows 93-98 down to 96-101.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert three new blank rows starting at row 93 (shifts old 93..98 -> 96..101)
$ws.Rows.Item(93).EntireRow.Insert()
$ws.Rows.Item(93).EntireRow.Insert()
$ws.Rows.Item(93).EntireRow.Insert()

# --- New row 93: Sandia, Extra, Región de O'Higgins ---
$ws.Cells.Item(93,1).Value  = 8
$ws.Cells.Item(93,2).Value  = "Terminal La Palmera de La Serena"
$ws.Cells.Item(93,3).Value  = "Coquimbo"
$ws.Cells.Item(93,4).Value  = 44610
$ws.Cells.Item(93,5).Value  = 4
$ws.Cells.Item(93,6).Value  = 100112028
$ws.Cells.Item(93,7).Value  = "Sandia"
$ws.Cells.Item(93,8).Value  = "Sin especificar"
$ws.Cells.Item(93,9).Value  = "Extra"
$ws.Cells.Item(93,10).Value = 3000
$ws.Cells.Item(93,11).Value = 3300
$ws.Cells.Item(93,12).Value = 3500
$ws.Cells.Item(93,13).Value = 3400
$ws.Cells.Item(93,14).Value = "`$/unidad"
$ws.Cells.Item(93,15).Value = "Región de O'Higgins"
$ws.Cells.Item(93,16).Value = 3400
$ws.Cells.Item(93,17).Value = 1
$ws.Cells.Item(93,18).Value = "Hortaliza"

# --- New row 94: Sandia, Primera, Región de O'Higgins ---
$ws.Cells.Item(94,1).Value  = 8
$ws.Cells.Item(94,2).Value  = "Terminal La Palmera de La Serena"
$ws.Cells.Item(94,3).Value  = "Coquimbo"
$ws.Cells.Item(94,4).Value  = 44610
$ws.Cells.Item(94,5).Value  = 4
$ws.Cells.Item(94,6).Value  = 100112028
$ws.Cells.Item(94,7).Value  = "Sandia"
$ws.Cells.Item(94,8).Value  = "Sin especificar"
$ws.Cells.Item(94,9).Value  = "Primera"
$ws.Cells.Item(94,10).Value = 3000
$ws.Cells.Item(94,11).Value = 2800
$ws.Cells.Item(94,12).Value = 3000
$ws.Cells.Item(94,13).Value = 2900
$ws.Cells.Item(94,14).Value = "`$/unidad"
$ws.Cells.Item(94,15).Value = "Región de O'Higgins"
$ws.Cells.Item(94,16).Value = 2900
$ws.Cells.Item(94,17).Value = 1
$ws.Cells.Item(94,18).Value = "Hortaliza"

# --- New row 95: Sandia, Segunda, Región de O'Higgins ---
$ws.Cells.Item(95,1).Value  = 8
$ws.Cells.Item(95,2).Value  = "Terminal La Palmera de La Serena"
$ws.Cells.Item(95,3).Value  = "Coquimbo"
$ws.Cells.Item(95,4).Value  = 44610
$ws.Cells.Item(95,5).Value  = 4
$ws.Cells.Item(95,6).Value  = 100112028
$ws.Cells.Item(95,7).Value  = "Sandia"
$ws.Cells.Item(95,8).Value  = "Sin especificar"
$ws.Cells.Item(95,9).Value  = "Segunda"
$ws.Cells.Item(95,10).Value = 2000
$ws.Cells.Item(95,11).Value = 2300
$ws.Cells.Item(95,12).Value = 2500
$ws.Cells.Item(95,13).Value = 2400
$ws.Cells.Item(95,14).Value = "`$/unidad"
$ws.Cells.Item(95,15).Value = "Región de O'Higgins"
$ws.Cells.Item(95,16).Value = 2400
$ws.Cells.Item(95,17).Value = 1
$ws.Cells.Item(95,18).Value = "Hortaliza"
